# Update TPM-derived NATMI metrics in the Wnt10a-Fzd8 LR-pairs sheet.
# Commit message: "update scripts wuth new tpm" -- ligand/receptor expression
# values (and everything derived from them) were recomputed with new TPM
# input, so refresh the affected cells in rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MuSCs -> MuSCs)
$ws.Range("G2").Value = 0.05606233333333333
$ws.Range("H2").Value = 0.168187
$ws.Range("M2").Value = 3.390429
$ws.Range("N2").Value = 10.171287
$ws.Range("O2").Value = 0.173121426386348
$ws.Range("P2").Value = 0.173121426386348
$ws.Range("Q2").Value = 0.190075360741
$ws.Range("R2").Value = 1.710678246669
$ws.Range("S2").Value = 0.173121426386348
$ws.Range("T2").Value = 0.173121426386348

# Row 3 (MuSCs -> ECs)
$ws.Range("G3").Value = 0.05606233333333333
$ws.Range("H3").Value = 0.168187
$ws.Range("O3").Value = 0.5936336753560868
$ws.Range("P3").Value = 0.5936336753560868
$ws.Range("Q3").Value = 0.6517687460563334
$ws.Range("R3").Value = 5.865918714507
$ws.Range("S3").Value = 0.5936336753560868
$ws.Range("T3").Value = 0.5936336753560868

# Row 4 (MuSCs -> FAPs)
$ws.Range("G4").Value = 0.05606233333333333
$ws.Range("H4").Value = 0.168187
$ws.Range("M4").Value = 4.546141666666667
$ws.Range("N4").Value = 13.638425
$ws.Range("O4").Value = 0.2321342018628743
$ws.Range("P4").Value = 0.2321342018628743
$ws.Range("Q4").Value = 0.2548673094972222
$ws.Range("R4").Value = 2.293805785475
$ws.Range("S4").Value = 0.2321342018628743
$ws.Range("T4").Value = 0.2321342018628743

# Row 5 (MuSCs -> Resolving-Mac)
$ws.Range("G5").Value = 0.05606233333333333
$ws.Range("H5").Value = 0.168187
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.021752
$ws.Range("N5").Value = 0.06525600000000001
$ws.Range("O5").Value = 0.001110696394691009
$ws.Range("P5").Value = 0.001110696394691009
$ws.Range("Q5").Value = 0.001219467874666667
$ws.Range("R5").Value = 0.010975210872
$ws.Range("S5").Value = 0.001110696394691009
$ws.Range("T5").Value = 0.001110696394691009
